$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "27.965.08"
$ws.Range("E2").Value = "  +0.18%  "

Set-TextValue "D3" "1.635.08"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.25%  "

Set-TextValue "D5" "212.14"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("E6").Value = "  -0.36%  "

Set-TextValue "D7" "0.997"
$ws.Range("E7").Value = "  -0.27%  "

Set-TextValue "D8" "23.43"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("E10").Value = "  -0.42%  "

Set-TextValue "D11" "0.0882"
$ws.Range("E11").Value = "  +0.65%  "

Set-TextValue "D12" "1.865.38"
$ws.Range("E12").Value = "  -0.64%  "

Set-TextValue "D13" "1.628.03"
$ws.Range("E13").Value = "  -0.99%  "

Set-TextValue "D14" "4.06"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("E15").Value = "  -2.32%  "

Set-TextValue "D16" "65.61"
$ws.Range("E16").Value = "  -0.37%  "

Set-TextValue "D17" "27.956.24"
$ws.Range("E17").Value = "  +0.12%  "

Set-TextValue "D18" "232.47"
$ws.Range("E18").Value = "  +1.01%  "

Set-TextValue "D19" "0.0₃0726"
$ws.Range("E19").Value = "  +0.06%  "

Set-TextValue "D20" "7.56"
$ws.Range("E20").Value = "  -1.01%  "

Set-TextValue "D21" "0.998"
$ws.Range("E21").Value = "  -0.28%  "

Set-TextValue "D22" "10.40"
$ws.Range("E22").Value = "  -4.58%  "

$ws.Range("E23").Value = "  -0.98%  "

$ws.Range("E24").Value = "  -3.10%  "

Set-TextValue "D25" "154.57"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -0.73%  "

Set-TextValue "D28" "15.66"
$ws.Range("E28").Value = "  -0.33%  "

Set-TextValue "D29" "0.998"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("E31").Value = "  -0.64%  "

Set-TextValue "D32" "3.41"
$ws.Range("E32").Value = "  +2.14%  "

Set-TextValue "D33" "1.412.11"
$ws.Range("E33").Value = "  -0.87%  "

$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  +8.79%  "

$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("E39").Value = "  +0.04%  "

Set-TextValue "D40" "0.872"
$ws.Range("E40").Value = "  -1.53%  "

$ws.Range("E41").Value = "  -1.14%  "

Set-TextValue "D42" "0.998"
$ws.Range("E42").Value = "  -0.23%  "

Set-TextValue "D43" "67.17"
$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("E44").Value = "  +1.16%  "

Set-TextValue "D45" "5.47"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("E46").Value = "  -0.36%  "

Set-TextValue "D47" "1.775.35"
$ws.Range("E47").Value = "  -0.58%  "

Set-TextValue "D48" "88.17"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D49" "0.0999"
$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.0505"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.54"
$ws.Range("E51").Value = "  -1.84%  "
